$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "001-16124695"
$ws.Range("B2").Value = "T16600471"
$ws.Range("C2").Value = "T16600471"

# Row 3
$ws.Range("A3").Value = "001-59566426"
$ws.Range("B3").Value = "23N0032641"
$ws.Range("C3").Value = "23N0032641"

# Row 4
$ws.Range("A4").Value = "001-55655832"
$ws.Range("B4").Value = "25M0085133"
$ws.Range("C4").Value = "25M0085133"

# Row 5
$ws.Range("A5").Value = "001-59549674"
$ws.Range("B5").Value = 2220563205
$ws.Range("C5").Value = 2220563205

# Row 6
$ws.Range("A6").Value = "001-59718794"
$ws.Range("B6").Value = "DJRDUA4279896"
$ws.Range("C6").Value = "25M0085202"

# Row 7
$ws.Range("A7").Value = "001-20584535"
$ws.Range("B7").Value = "DJAUSA4279777"
$ws.Range("C7").Value = "24N0016923"

# Row 8
$ws.Range("A8").Value = "001-59221260"
$ws.Range("B8").Value = "DJLRDA4279757"
$ws.Range("C8").Value = 211084425

# Row 9
$ws.Range("A9").Value = "001-20313274"
$ws.Range("B9").Value = "DJPHXA4279754"
$ws.Range("C9").Value = 2110505708

# Row 10
$ws.Range("A10").Value = "001-59521162"
$ws.Range("B10").Value = "DJAUSA4279728"
$ws.Range("C10").Value = "24N0016913"

# Row 11
$ws.Range("A11").Value = "001-59851956"
$ws.Range("B11").Value = 2110505575
$ws.Range("C11").Value = 2110505575
